# Apply the "DeudoresPrueba" update: refresh the client/payment rows with
# the new data set (32 data rows instead of 38) and drop the trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the six rows that no longer exist in the updated data (old rows
# 34-39); everything below shifts up and the sheet's used range becomes
# A1:E33.
$ws.Rows("34:39").Delete()

# Rewrite the 32 data rows (row 1 is the header and is unchanged).
$data = @(
    @(1,  "ALISO",              45996, 196000),
    @(2,  "CAMILIN",             45997, 166000),
    @(3,  "CARNILANDIA",         46000, 436000),
    @(4,  "CARNIVOROS",          45959, 200000),
    @(5,  "CIMARRON DORADO",     46000, 473300),
    @(6,  "CLIENTE PAOLA",       46000, 92000),
    @(7,  "COCINA CHINA",        45998, 170000),
    @(8,  "DARWIN FUTBOL",       45921, 200000),
    @(9,  "DAVIDCITO",           45947, 100000),
    @(10, "EL RES",              45997, 240000),
    @(11, "EL RUBY",             45992, 85100),
    @(12, "FRANCO",              45996, 20000),
    @(13, "LA SELECTA",          45912, 82000),
    @(14, "MERKA FRUVER DEXI",   45995, 339000),
    @(15, "MERKA FRUVER DEXI",   45988, 115400),
    @(16, "NEVADA",              46000, 164000),
    @(17, "NEVADA",              45996, 229000),
    @(18, "NOVILLON SAN MATEO",  45971, 83000),
    @(19, "PARAÍSO FUNZA",       45996, 202000),
    @(20, "PARAÍSO MOSQUERA",    46000, 354900),
    @(21, "PINILLA",             45924, 16000),
    @(22, "PINILLA",             45931, 166000),
    @(23, "PINILLA SOACHA",      45993, 129000),
    @(24, "PLAZA JESSICA",       45999, 971300),
    @(25, "PLAZA JESSICA",       46000, 1238000),
    @(26, "PORTAL ZIPA",         45995, 224000),
    @(27, "PUNTA DE ANCA",       46000, 7600),
    @(28, "SANDRA 20 DE JULIO",  46000, 300000),
    @(29, "SANTANDER MADRID",    45996, 63000),
    @(30, "SANTANDER SUR",       45997, 250700),
    @(31, "SANTANDER SUR",       45993, 80000),
    @(32, "VNZLNO PUNTA ANCA",   45992, 82000)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $item = $data[$i]

    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $false
}

Write-Output "DeudoresPrueba rows refreshed"
